# AP z boson CMS 8 TeV
# Adds a new "process" column (Y) to sheet1, with header "process" in Y1
# and value "pp->Z/gamma*->l+ l-" for every data row (Y2:Y35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Y1").Value = "process"

$lastRow = 35
$ws.Range("Y2:Y$lastRow").Value = "pp->Z/gamma*->l+ l-"

[void]$ws.Range("R1").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 18
[void]$ws.Range("AB35").Select()
